$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 30.4397731227338
    3 = 38.20502188407325
    4 = 35.21436270867769
    5 = 48.43409545658513
    6 = 61.07714571501362
    7 = 70.50510813979032
    8 = 80.53375558817824
    9 = 87.52339934547867
    10 = 94.162789435116
    11 = 109.9088414734317
    12 = 115.3182243824334
    13 = 126.7658016717521
    14 = 136.2272366785352
    15 = 146.0478662119301
    16 = 153.5387168450899
    17 = 166.5430047247741
    18 = 173.9547299083377
    19 = 190.3025432890343
    20 = 200.4110385373977
    21 = 204.8282086163355
    22 = 212.6839822225405
    23 = 219.7050737278354
    24 = 244.0190363075733
    25 = 247.6221600294006
    26 = 258.3740680441053
    27 = 268.0018304445429
    28 = 278.2614259565154
    29 = 283.7885154259984
    30 = 296.4583497063255
    31 = 309.5153359641639
    32 = 304.3426334569814
    33 = 307.3211362957871
    34 = 324.3294911169081
    35 = 339.130649220143
    36 = 339.6680929850742
    37 = 353.9380312409988
    38 = 364.0731227609003
    39 = 376.6321761199234
    40 = 385.1840067871106
    41 = 402.1141205234409
    42 = 409.6757775777224
    43 = 417.7259465236669
    44 = 427.9950731796807
    45 = 442.303775585897
    46 = 451.2115804612736
    47 = 464.1253751266472
    48 = 462.7684141179541
    49 = 532.4895523720602
    50 = 538.5350260499763
    51 = 546.4948659263106
    52 = 557.3356022650211
    53 = 561.0766281123863
    54 = 570.5794059401441
    55 = 577.8494618276045
    56 = 583.6887118733938
    57 = 592.8609074834558
    58 = 590.8992244114816
    59 = 601.494854113642
    60 = 601.5502142789176
    61 = 601.535985307298
    62 = 601.5631457759213
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
}
